# chore: adapt column header formatting to respective input file names
#
# The header row (row 1) uses two parallel sets of column names: the first
# block (columns A-J) is suffixed "_old" and describes the "before" input
# file, the second block (columns L-U) is suffixed "_new" and describes the
# "after" input file. Column K just holds the literal header "diff".
#
# Rename the suffixes to the concrete format-version identifiers that the
# two input files actually correspond to: "_old" -> "_FV2210" and
# "_new" -> "_FV2304". Then turn the sheet into a proper Excel Table
# (ListObject) and freeze the header row so it stays visible while
# scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{ "_old" = "_FV2210"; "_new" = "_FV2304" }

$usedRange = $ws.UsedRange
$firstCol = $usedRange.Column
$lastCol = $firstCol + $usedRange.Columns.Count - 1

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Value2
    if ($null -eq $header) { continue }

    foreach ($suffix in $suffixMap.Keys) {
        if ($header.EndsWith($suffix)) {
            $base = $header.Substring(0, $header.Length - $suffix.Length)
            $cell.Value2 = $base + $suffixMap[$suffix]
            break
        }
    }
}

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table so headers get filter buttons and
# consistent formatting.
$dataRange = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"
